$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.422.72"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").Value = "3.692.88"
$ws.Range("E3").Value = "  -2.99%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "691.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.49%  "
$ws.Range("D7").Value = "3.691.33"
$ws.Range("E7").Value = "  -3.03%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -4.79%  "
$ws.Range("E10").Value = "  -8.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.17%  "
$ws.Range("E13").Value = "  -5.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.15%  "
$ws.Range("D15").Value = "4.315.72"
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").Value = "3.687.26"
$ws.Range("E16").Value = "  -3.25%  "
$ws.Range("D17").Value = "69.429.61"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "480.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.38%  "
$ws.Range("E22").Value = "  -5.79%  "
$ws.Range("E24").Value = "  -4.80%  "
$ws.Range("D25").Value = "3.839.47"
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("E26").Value = "  -9.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.63%  "
$ws.Range("E31").Value = "  -10.11%  "
$ws.Range("E32").Value = "  -7.66%  "
$ws.Range("E33").Value = "  -7.60%  "
$ws.Range("E34").Value = "  -5.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("E36").Value = "  -7.06%  "
$ws.Range("D37").Value = "3.660.70"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("E38").Value = "  -7.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("E41").Value = "  -7.94%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  -6.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "30.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000284"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.94%  "
